$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 232, pushing existing rows 232-332 down to 233-333.
$ws.Rows.Item(232).Insert()

# Columns A,B,C,E,F,G,H,I,N,O,Q,R are identical between the old row 232 (now shifted
# to row 233) and the new row that needs to be created at row 232, so copy them across.
$colsToCopy = @(1,2,3,5,6,7,8,9,14,15,17,18)
foreach ($col in $colsToCopy) {
    $ws.Cells.Item(232, $col).Value2 = $ws.Cells.Item(233, $col).Value2
}

# Match the date cell's number format/style used by the rest of column D.
$ws.Cells.Item(232, 4).NumberFormat = $ws.Cells.Item(233, 4).NumberFormat

# New record values for the inserted row (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg).
$ws.Cells.Item(232, 4).Value2 = 44845
$ws.Cells.Item(232, 10).Value2 = 400
$ws.Cells.Item(232, 11).Value2 = 31000
$ws.Cells.Item(232, 12).Value2 = 31000
$ws.Cells.Item(232, 13).Value2 = 31000
$ws.Cells.Item(232, 16).Value2 = 517
